$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1727
$ws.Range("I19").Value = 2122.5
$ws.Range("J19").Value = 1199.6666
$ws.Range("K19").Value = 2122.5
$ws.Range("L19").Value = 1199.6666
$ws.Range("M19").Value = -1947.5
$ws.Range("N19").Value = -1549.6666

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 20771.4
$ws.Range("I96").Value = 34176
$ws.Range("K96").Value = 102528
$ws.Range("M96").Value = -101155

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3020.2222
$ws.Range("I132").Value = 1624.4667
$ws.Range("K132").Value = 4873.4001
$ws.Range("M132").Value = -2343.4001

# ARM row 3
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 443.66666
$ws.Range("I5").Value = 492.6
$ws.Range("K5").Value = 492.6
$ws.Range("M5").Value = -380.6

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4242.778

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1137
$ws.Range("I74").Value = 1164.4
$ws.Range("K74").Value = 1164.4
$ws.Range("M74").Value = -290.4000000000001

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1137
$ws.Range("I77").Value = 1164.4
$ws.Range("K77").Value = 5822
$ws.Range("M77").Value = -1454

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5399.4
$ws.Range("I102").Value = 5249.25
$ws.Range("K102").Value = 5249.25
$ws.Range("M102").Value = -3627.25

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3521.5833
$ws.Range("I110").Value = 3659.9092
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 3659.9092
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = -1614.9092
$ws.Range("N110").Value = -6090

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2245.6365
$ws.Range("I122").Value = 1359.6
$ws.Range("K122").Value = 4078.8
$ws.Range("M122").Value = -1628.8

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3848.5386
$ws.Range("J132").Value = 3334.6667
$ws.Range("L132").Value = 10004.0001
$ws.Range("N132").Value = -15064.0001

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 443.66666
$ws.Range("I4").Value = 492.6
$ws.Range("K4").Value = 492.6
$ws.Range("M4").Value = -377.6

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3624.25
$ws.Range("I86").Value = 3624.25
$ws.Range("K86").Value = 3624.25
$ws.Range("M86").Value = -2501.25

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3624.25
$ws.Range("I89").Value = 3624.25
$ws.Range("K89").Value = 18121.25
$ws.Range("M89").Value = -12505.25

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4699.077
$ws.Range("I107").Value = 4694.364
$ws.Range("J107").Value = 4725
$ws.Range("K107").Value = 4694.364
$ws.Range("L107").Value = 4725
$ws.Range("M107").Value = -2774.364
$ws.Range("N107").Value = -8565

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 77890
$ws.Range("J140").Value = 77890
$ws.Range("L140").Value = 77890
$ws.Range("N140").Value = -88250

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 109
$ws.Range("I7").Value = 68.333336
$ws.Range("J7").Value = 157.8
$ws.Range("K7").Value = 68.333336
$ws.Range("L7").Value = 157.8
$ws.Range("M7").Value = 44.666664
$ws.Range("N7").Value = -383.8

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 534.2727
$ws.Range("J22").Value = 679.75
$ws.Range("L22").Value = 679.75
$ws.Range("N22").Value = -1379.75

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3051.5833
$ws.Range("I31").Value = 3888.5715
$ws.Range("J31").Value = 1879.8
$ws.Range("K31").Value = 3888.5715
$ws.Range("L31").Value = 1879.8
$ws.Range("M31").Value = -3593.5715
$ws.Range("N31").Value = -2469.8

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3051.5833
$ws.Range("I34").Value = 3888.5715
$ws.Range("J34").Value = 1879.8
$ws.Range("K34").Value = 3888.5715
$ws.Range("L34").Value = 1879.8
$ws.Range("M34").Value = -3686.5715
$ws.Range("N34").Value = -2283.8

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1569.5
$ws.Range("I134").Value = 1356.5
$ws.Range("K134").Value = 4069.5
$ws.Range("M134").Value = -1534.5

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 69951.336
$ws.Range("J37").Value = 69951.336
$ws.Range("L37").Value = 209854.008
$ws.Range("N37").Value = -210078.008

# CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 5666
$ws.Range("J75").Value = 8003
$ws.Range("L75").Value = 24009
$ws.Range("N75").Value = -26005

# CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 5666
$ws.Range("J78").Value = 8003
$ws.Range("L78").Value = 72027
$ws.Range("N78").Value = -82011

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 14250
$ws.Range("J5").Value = 15666.667
$ws.Range("L5").Value = 15666.667
$ws.Range("N5").Value = -15890.667

# GSM row 12
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 15001.111
$ws.Range("I12").Value = 15001.5
$ws.Range("J12").Value = 15000.8
$ws.Range("K12").Value = 15001.5
$ws.Range("L12").Value = 15000.8
$ws.Range("M12").Value = -14861.5
$ws.Range("N12").Value = -15280.8

# GSM row 18
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 36672
$ws.Range("I18").Value = 30005
$ws.Range("J18").Value = 50006
$ws.Range("K18").Value = 30005
$ws.Range("L18").Value = 50006
$ws.Range("M18").Value = -29712
$ws.Range("N18").Value = -50592

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7666.3335
$ws.Range("I70").Value = 7500
$ws.Range("K70").Value = 7500
$ws.Range("M70").Value = -7230

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7666.3335
$ws.Range("I73").Value = 7500
$ws.Range("K73").Value = 7500
$ws.Range("M73").Value = -6564

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4500
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 5000
$ws.Range("N80").Value = -6996

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4500
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 25000
$ws.Range("N83").Value = -34984

# GSM row 94
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 33990.43
$ws.Range("J94").Value = 33990.43
$ws.Range("L94").Value = 33990.43
$ws.Range("N94").Value = -35342.43

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 990
$ws.Range("I113").Value = 990
$ws.Range("K113").Value = 990
$ws.Range("M113").Value = 1180

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3387.2
$ws.Range("I122").Value = 3369.077
$ws.Range("K122").Value = 10107.231
$ws.Range("M122").Value = -7657.231

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4075.6667
$ws.Range("I132").Value = 3290.8
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 9872.400000000001
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -7342.400000000001
$ws.Range("N132").Value = -29060

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5250.6665
$ws.Range("I61").Value = 5002
$ws.Range("J61").Value = 5375
$ws.Range("K61").Value = 5002
$ws.Range("L61").Value = 5375
$ws.Range("M61").Value = -4800
$ws.Range("N61").Value = -5779

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5250.6665
$ws.Range("I113").Value = 5002
$ws.Range("J113").Value = 5375
$ws.Range("K113").Value = 5002
$ws.Range("L113").Value = 5375
$ws.Range("M113").Value = -2832
$ws.Range("N113").Value = -9715

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7338.515
$ws.Range("I122").Value = 7149.6113
$ws.Range("J122").Value = 7565.2
$ws.Range("K122").Value = 21448.8339
$ws.Range("L122").Value = 22695.6
$ws.Range("M122").Value = -18998.8339
$ws.Range("N122").Value = -27595.6

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3895.348
$ws.Range("I132").Value = 4234.6665
$ws.Range("J132").Value = 2673.8
$ws.Range("K132").Value = 12703.9995
$ws.Range("L132").Value = 8021.400000000001
$ws.Range("M132").Value = -10173.9995
$ws.Range("N132").Value = -13081.4

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2257.7
$ws.Range("I136").Value = 2225.2222
$ws.Range("K136").Value = 6675.6666
$ws.Range("M136").Value = -4125.6666

# WVR row 70
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 87792.5
$ws.Range("I70").Value = 87795
$ws.Range("J70").Value = 87790
$ws.Range("K70").Value = 87795
$ws.Range("L70").Value = 87790
$ws.Range("M70").Value = -87480
$ws.Range("N70").Value = -88420

# WVR row 73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 87792.5
$ws.Range("I73").Value = 87795
$ws.Range("J73").Value = 87790
$ws.Range("K73").Value = 87795
$ws.Range("L73").Value = 87790
$ws.Range("M73").Value = -86703
$ws.Range("N73").Value = -89974

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6208.9473
$ws.Range("I132").Value = 6442.778
$ws.Range("K132").Value = 19328.334
$ws.Range("M132").Value = -16798.334

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5314.56
$ws.Range("I136").Value = 6040.579
$ws.Range("J136").Value = 3015.5
$ws.Range("K136").Value = 18121.737
$ws.Range("L136").Value = 9046.5
$ws.Range("M136").Value = -15571.737
$ws.Range("N136").Value = -14146.5
